$wb = $excel.ActiveWorkbook

$wsShortTerm = $wb.Worksheets.Item("short term")
$wsDone      = $wb.Worksheets.Item("done")

# 1. Fix the wording of task 86 (still sitting in "short term", row 21)
$wsShortTerm.Range("A21").Value = "86. Add on the Dromics web page an introduction of the functions to help biological interpretation (and sva intsallation - no need just suggested )"

# 2. The task is now finished: move the whole row (task text + owner) from
#    "short term" (row 21) down to the bottom of the "done" sheet (row 66),
#    keeping its border/wrap formatting but dropping the "short term" highlight fill.
$src = $wsShortTerm.Range("A21:B21")
$dst = $wsDone.Range("A66:B66")

$src.Copy($dst)
$dst.Interior.ColorIndex = -4142

# 3. Remove the now-duplicated row from "short term" - everything below shifts up.
$wsShortTerm.Rows.Item(21).Delete()

# Leave the selections where the author would naturally have left them:
# the freshly dropped row in "done", then back to the shifted-up row in "short term".
$wsDone.Range("A66:B66").Select()
$wsShortTerm.Rows.Item(21).Select()
